# Add data for 2021-10-19
# Rename the sheet and update the header text to reflect the new "through" date,
# then update/insert the individual neighborhood/month count cells per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab
$ws.Name = "Through 2021-10-11"

# Update the title cell text (shared string used in A1's column header row... actually
# it is the first shared string, used as a column header in row 1)
$ws.Range("B1").Value = "October 2021 (through October 11)"

# Row 2 (Garfield Park)
$ws.Range("B2").Value = 7
$ws.Range("BJ2").Value = 3

# Row 4 (Chatham)
$ws.Range("B4").Value = 6
$ws.Range("AF4").Value = 1

# Row 8 (Grand Crossing)
$ws.Range("AF8").Value = 1

# Row 14 (Little Italy, UIC)
$ws.Range("L14").Value = 1

# Row 15 (Lake View)
$ws.Range("B15").Value = 5

# Row 21 (United Center)
$ws.Range("L21").Value = 2

# Row 24 (Avalon Park)
$ws.Range("B24").Value = 1

# Row 31 (Ashburn)
$ws.Range("L31").Value = 2
$ws.Range("AF31").Value = 3

# Row 37 (Rogers Park)
$ws.Range("AZ37").Value = 3

# Row 40 (Pullman)
$ws.Range("AP40").Value = 1

# Row 42 (West Loop)
$ws.Range("L42").Value = 1

# Row 48 (Archer Heights)
$ws.Range("L48").Value = 4

# Row 57 (New City)
$ws.Range("B57").Value = 1

# Row 77 (Jefferson Park)
$ws.Range("L77").Value = 1
